# Galot Data Log - "final raw data from Galot"
# Adds one more "SD CARD SWAP" section (row 37) followed by its Run 0/Run 1/Run 2
# spindown notes (rows 38-41), extending the data log table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37: new "SD CARD SWAP" section header (bold, matching the other section headers)
$ws.Range("B37").Value = "SD CARD SWAP"
$ws.Range("B37").Font.Bold = $true

# Row 38: Run 0
$ws.Range("B38").Value = "Run 0"
$ws.Range("C38").Value = "spindown @ 10m/s"

# Row 39: Run 1
$ws.Range("B39").Value = "Run 1"
$ws.Range("C39").Value = "spindown @ 10m/s"
$ws.Range("D39").Value = "downhill"

# back to filling in the rest of row 38
$ws.Range("D38").Value = "uphill"
$ws.Range("E38").Value = "along the wider section of the road"
$ws.Range("H38").Value = "spindowns5.TXT"

# Row 40: Run 2
$ws.Range("B40").Value = "Run 2"
$ws.Range("C40").Value = "spindown on sticky"
$ws.Range("D40").Value = "uphill"

# Row 41: continuation of Run 2 (tire popped at the end)
$ws.Range("C41").Value = "spindown on sticky"
$ws.Range("D41").Value = "downhill"
$ws.Range("E41").Value = "TIRE POPPED AT END"

# Move the selection/active cell and scroll position to the newly entered data
$ws.Range("E41").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
